# Add data for 2024-03-10
# Updates the year-to-date crime-count cells across the citywide, by-neighborhood
# and per-neighborhood sheets to reflect one additional day of data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("B2").Value = 8
$ws.Range("D2").Value = 12
$ws.Range("F2").Value = 14
$ws.Range("I2").Value = 12
$ws.Range("K2").Value = 22
$ws.Range("I3").Value = 33
$ws.Range("J3").Value = 31
$ws.Range("E6").Value = 1
$ws.Range("C9").Value = 83
$ws.Range("E9").Value = 82
$ws.Range("F9").Value = 110
$ws.Range("I9").Value = 95
$ws.Range("J9").Value = 79
$ws.Range("K9").Value = 72
$ws.Range("B10").Value = 174
$ws.Range("C10").Value = 204
$ws.Range("D10").Value = 339
$ws.Range("E10").Value = 327
$ws.Range("F10").Value = 474
$ws.Range("G10").Value = 367
$ws.Range("I10").Value = 141
$ws.Range("J10").Value = 108
$ws.Range("K10").Value = 133
$ws.Range("B11").Value = 282
$ws.Range("C11").Value = 317
$ws.Range("D11").Value = 462
$ws.Range("E11").Value = 444
$ws.Range("F11").Value = 614
$ws.Range("G11").Value = 515
$ws.Range("I11").Value = 287
$ws.Range("J11").Value = 253
$ws.Range("K11").Value = 263

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("D2").Value = 2
$ws.Range("G5").Value = 18
$ws.Range("K5").Value = 4
$ws.Range("D6").Value = 17
$ws.Range("G6").Value = 30
$ws.Range("K6").Value = 13

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K6").Value = 5
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 11
$ws.Range("K8").Value = 14

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("F2").Value = 2
$ws.Range("K2").Value = 2
$ws.Range("F5").Value = 15
$ws.Range("J5").Value = 15
$ws.Range("K5").Value = 7
$ws.Range("B6").Value = 18
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 89
$ws.Range("G6").Value = 70
$ws.Range("I6").Value = 39
$ws.Range("K6").Value = 20
$ws.Range("B7").Value = 27
$ws.Range("D7").Value = 116
$ws.Range("E7").Value = 106
$ws.Range("F7").Value = 156
$ws.Range("G7").Value = 88
$ws.Range("I7").Value = 66
$ws.Range("J7").Value = 41
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 10

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("B7").Value = 9
$ws.Range("G7").Value = 16
$ws.Range("E8").Value = 18
$ws.Range("B27").Value = 19
$ws.Range("C27").Value = 24
$ws.Range("E27").Value = 12
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 11
$ws.Range("D31").Value = 17
$ws.Range("G31").Value = 30
$ws.Range("K31").Value = 13
$ws.Range("B35").Value = 11
$ws.Range("K35").Value = 14
$ws.Range("E46").Value = 11
$ws.Range("F46").Value = 10
$ws.Range("I46").Value = 3
$ws.Range("J46").Value = 7
$ws.Range("E49").Value = 10
$ws.Range("B52").Value = 27
$ws.Range("D52").Value = 116
$ws.Range("E52").Value = 106
$ws.Range("F52").Value = 156
$ws.Range("G52").Value = 88
$ws.Range("I52").Value = 66
$ws.Range("J52").Value = 41
$ws.Range("K52").Value = 37
$ws.Range("D67").Value = 7
$ws.Range("F67").Value = 9
$ws.Range("E76").Value = 11
$ws.Range("I76").Value = 11
$ws.Range("C77").Value = 3
$ws.Range("E77").Value = 9
$ws.Range("F78").Value = 7
$ws.Range("D87").Value = 4
$ws.Range("E90").Value = 10
$ws.Range("G90").Value = 9
$ws.Range("F93").Value = 19
$ws.Range("E94").Value = 5
$ws.Range("B97").Value = 282
$ws.Range("C97").Value = 317
$ws.Range("D97").Value = 462
$ws.Range("E97").Value = 444
$ws.Range("F97").Value = 614
$ws.Range("G97").Value = 515
$ws.Range("I97").Value = 287
$ws.Range("J97").Value = 253
$ws.Range("K97").Value = 263

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 4

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 7

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("C3").Value = 3
$ws.Range("E4").Value = 9
$ws.Range("C5").Value = 3
$ws.Range("E5").Value = 9

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 3
$ws.Range("B5").Value = 14
$ws.Range("C5").Value = 17
$ws.Range("E5").Value = 9
$ws.Range("J5").Value = 4
$ws.Range("B6").Value = 19
$ws.Range("C6").Value = 24
$ws.Range("E6").Value = 12
$ws.Range("I6").Value = 14
$ws.Range("J6").Value = 11

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J3").Value = 3
$ws.Range("E5").Value = 1
$ws.Range("I5").Value = 2
$ws.Range("F6").Value = 7
$ws.Range("E7").Value = 11
$ws.Range("F7").Value = 10
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 7

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("E4").Value = 7
$ws.Range("G4").Value = 9
$ws.Range("E5").Value = 10
$ws.Range("G5").Value = 9

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 1
$ws.Range("E7").Value = 7
$ws.Range("E8").Value = 11
$ws.Range("I8").Value = 11

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("F5").Value = 17
$ws.Range("F6").Value = 19

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 5

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("B2").Value = 1
$ws.Range("G5").Value = 13
$ws.Range("B6").Value = 9
$ws.Range("G6").Value = 16

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("D6").Value = 7
$ws.Range("F6").Value = 8
$ws.Range("D7").Value = 7
$ws.Range("F7").Value = 9

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("E4").Value = 1
$ws.Range("E7").Value = 18
